# TiffinTracker.xlsx edit
# Commit message: "Add files via upload"
#
# Semantic changes captured by the diff:
#  - Row 15 (Friday, 2018-01-12) tiffin flags for Tejashri, Shradha, Rahul P,
#    Rahul G, Vimarsh, Taniya, Divyam and Siddesh (columns D:K) are filled
#    in with 1, and the Notes-adjacent "Soumya" column (M) is also set to 1.
#    Column L (Lovely) is left untouched (it already holds 0).
#    The downstream SUM()/totals formulas on rows 35-39 recompute
#    automatically from these new values.
#  - The active selection on the January sheet moves to O16 (and the sheet
#    no longer has a frozen/scrolled topLeftCell override).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("January")

$ws.Range("D15:K15").Value = 1
$ws.Range("M15").Value = 1

$ws.Range("O16").Select() | Out-Null
